# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp
# - Refresh scraped COVID counters for several countries (rows re-rank by
#   total cases, so Pakistan/Banglades and Jamaica/Georgia trade places)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 07:57"

# Row 18/19: Pakistan overtakes Banglades
$ws.Range("A18").Value = "Pakistan"
$ws.Range("B18").Value2 = 292765
$ws.Range("C18").Value2 = 591
$ws.Range("D18").Value2 = 275836
$ws.Range("E18").Value2 = 10694
$ws.Range("F18").Value2 = 0
$ws.Range("G18").Value2 = 4
$ws.Range("H18").Value2 = 6235

$ws.Range("A19").Value = "Banglades"
$ws.Range("B19").Value2 = 292625
$ws.Range("C19").Value2 = 0
$ws.Range("D19").Value2 = 175567
$ws.Range("E19").Value2 = 113151
$ws.Range("F19").Value2 = 0
$ws.Range("G19").Value2 = 0
$ws.Range("H19").Value2 = 3907

# Row 33: Israel data refresh
$ws.Range("B33").Value2 = 102080
$ws.Range("C33").Value2 = 147
$ws.Range("D33").Value2 = 79303
$ws.Range("E33").Value2 = 21958

# Row 56: Kirguistan data refresh
$ws.Range("B56").Value2 = 43023
$ws.Range("C56").Value2 = 134
$ws.Range("D56").Value2 = 36397
$ws.Range("E56").Value2 = 5570
$ws.Range("G56").Value2 = 1
$ws.Range("H56").Value2 = 1056

# Row 62: Uzbekistan data refresh
$ws.Range("B62").Value2 = 38698
$ws.Range("C62").Value2 = 286
$ws.Range("E62").Value2 = 3853
$ws.Range("G62").Value2 = 4
$ws.Range("H62").Value2 = 269

# Row 120: Tailandia data refresh
$ws.Range("B120").Value2 = 3395
$ws.Range("C120").Value2 = 5
$ws.Range("D120").Value2 = 3221
$ws.Range("E120").Value2 = 116

# Row 142: Bahamas data refresh
$ws.Range("B142").Value2 = 1765
$ws.Range("C142").Value2 = 62
$ws.Range("D142").Value2 = 227
$ws.Range("E142").Value2 = 1509
$ws.Range("G142").Value2 = 2
$ws.Range("H142").Value2 = 29

# Row 149/150: Jamaica overtakes Georgia
$ws.Range("A149").Value = "Jamaica"
$ws.Range("B149").Value2 = 1413
$ws.Range("C149").Value2 = 67
$ws.Range("D149").Value2 = 817
$ws.Range("E149").Value2 = 580
$ws.Range("F149").Value2 = 0
$ws.Range("G149").Value2 = 0
$ws.Range("H149").Value2 = 16

$ws.Range("A150").Value = "Georgia"
$ws.Range("B150").Value2 = 1394
$ws.Range("C150").Value2 = 0
$ws.Range("D150").Value2 = 1132
$ws.Range("E150").Value2 = 245
$ws.Range("F150").Value2 = 0
$ws.Range("G150").Value2 = 0
$ws.Range("H150").Value2 = 17
